$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" '42.047.27'
Set-TextCell "E2" '  -2.34%  '
Set-TextCell "D3" '2.231.49'
Set-TextCell "E3" '  -3.14%  '
Set-TextCell "E4" '  +0.20%  '
Set-TextCell "D5" '246.01'
Set-TextCell "E5" '  -2.67%  '
Set-TextCell "D6" '0.631'
Set-TextCell "E6" '  -1.89%  '
Set-TextCell "D7" '76.32'
Set-TextCell "E7" '  +0.10%  '
Set-TextCell "D9" '0.618'
Set-TextCell "E9" '  -5.38%  '
Set-TextCell "D10" '41.37'
Set-TextCell "E10" '  +4.20%  '
Set-TextCell "D11" '0.0943'
Set-TextCell "E11" '  -4.63%  '
Set-TextCell "D12" '7.04'
Set-TextCell "E12" '  -9.34%  '
Set-TextCell "E13" '  -3.71%  '
Set-TextCell "D14" '2.571.20'
Set-TextCell "E14" '  -2.85%  '
Set-TextCell "D15" '14.64'
Set-TextCell "E15" '  -5.61%  '
Set-TextCell "E16" '  -3.50%  '
Set-TextCell "D17" '2.233.12'
Set-TextCell "E17" '  -3.22%  '
Set-TextCell "D18" '41.895.47'
Set-TextCell "E18" '  -2.57%  '
Set-TextCell "E19" '  -3.83%  '
Set-TextCell "D20" '71.50'
Set-TextCell "E20" '  -2.10%  '
Set-TextCell "D21" '6.06'
Set-TextCell "E21" '  -3.95%  '
Set-TextCell "D22" '2.28'
Set-TextCell "E22" '  +0.60%  '
Set-TextCell "D23" '230.34'
Set-TextCell "E23" '  -3.61%  '
Set-TextCell "E24" '  -0.18%  '
Set-TextCell "B25" 'WEMIXToken'
Set-TextCell "C25" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell "D25" '3.69'
Set-TextCell "E25" '  -5.55%  '
Set-TextCell "B26" 'Cosmos'
Set-TextCell "C26" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell "D26" '11.18'
Set-TextCell "E26" '  -4.03%  '
Set-TextCell "E27" '  -6.04%  '
Set-TextCell "D28" '7.37'
Set-TextCell "E28" '  +14.82%  '
Set-TextCell "E29" '  -1.99%  '
Set-TextCell "D30" '168.77'
Set-TextCell "E30" '  +0.48%  '
Set-TextCell "D31" '20.50'
Set-TextCell "E31" '  -3.69%  '
Set-TextCell "D32" '0.0824'
Set-TextCell "E32" '  -2.53%  '
Set-TextCell "D33" '32.25'
Set-TextCell "E33" '  +4.66%  '
Set-TextCell "E34" '  -6.99%  '
Set-TextCell "E35" '  -2.29%  '
Set-TextCell "D36" '4.46'
Set-TextCell "E36" '  -3.73%  '
Set-TextCell "D37" '4.93'
Set-TextCell "E37" '  +1.38%  '
Set-TextCell "E38" '  -4.58%  '
Set-TextCell "D39" '13.96'
Set-TextCell "E39" '  +0.44%  '
Set-TextCell "D40" '5.82'
Set-TextCell "E40" '  -1.20%  '
Set-TextCell "E41" '  -8.45%  '
Set-TextCell "D42" '112.43'
Set-TextCell "E42" '  +7.25%  '
Set-TextCell "E43" '  -8.36%  '
Set-TextCell "D44" '60.39'
Set-TextCell "E44" '  -3.85%  '
Set-TextCell "D45" '8.64'
Set-TextCell "E45" '  -6.50%  '
Set-TextCell "D46" '0.0987'
Set-TextCell "E46" '  -4.72%  '
Set-TextCell "D47" '0.996'
Set-TextCell "E47" '  -0.50%  '
Set-TextCell "E48" '  -4.55%  '
Set-TextCell "E49" '  -2.68%  '
Set-TextCell "B50" 'WOONetwork'
Set-TextCell "C50" 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextCell "D50" '0.434'
Set-TextCell "E50" '  +13.88%  '
Set-TextCell "B51" 'FTXToken'
Set-TextCell "C51" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell "D51" '4.17'
Set-TextCell "E51" '  -15.32%  '
